# Phase 2 implementation checklist: check off several checklist boxes by
# inserting "x" characters inside the "[ ]" markers.
#
# To match the target authoring pattern, each inserted "x" needs to land in
# its own run, split away from the surrounding literal text. This host (like
# Word itself) fuses adjacent runs back together at save time whenever they
# are adjacent AND share identical run formatting - so a plain InsertBefore
# merges the new character straight into whatever run it landed in. Briefly
# toggling a character-formatting property (Bold) on just the inserted
# character and then reverting it forces a permanent run boundary around it
# (the two sides no longer look "freshly typed into the same run") without
# changing its visible formatting.
#
# When a paragraph needs more than one inserted "x" (the "Create Phase 2
# working directory" line gets two), the formatting toggles are applied from
# right to left, after all the plain-text inserts are done. Toggling
# left-to-right would re-merge an already-split "x" run with its new
# not-yet-split neighbour to its right before that neighbour gets its own
# split; going right to left means every toggle's right-hand neighbour is
# already final by the time it runs.

$d = $word.ActiveDocument

function Find-Start($SearchText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($SearchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $SearchText"
        return -1
    }
    return $rng.Start
}

function Split-Char($Pos) {
    # Force a run boundary around the single character at $Pos by toggling
    # Bold on and back off.
    $xr = $d.Range($Pos, $Pos + 1)
    $xr.Font.Bold = 1
    $xr.Font.Bold = 0
}

function Insert-CheckMark($SearchText, $Offset) {
    # Single "x" insertion at $Offset chars into the match for $SearchText.
    $start = Find-Start $SearchText
    if ($start -lt 0) { return }
    $pos = $start + $Offset
    $d.Range($pos, $pos).InsertBefore("x")
    Split-Char $pos
}

# --- Pre-Implementation Checklist ---------------------------------------
# "[ ] Backup current enhanced_hospitals.yaml"  ->  "[ x] Backup current ..."
Insert-CheckMark "[ ] Backup current enhanced_hospitals.yaml" 2

# "[ ] Backup current pattern_based_scraper.R"  ->  "[ x] Backup current ..."
Insert-CheckMark "[ ] Backup current pattern_based_scraper.R" 2

# "[ ] Create baseline snapshot of all hospitals"  ->  "[ x] Create baseline ..."
Insert-CheckMark "[ ] Create baseline snapshot of all hospitals" 2

# "[ ] Create Phase 2 working directory"  ->  "[ x]x Create Phase 2 working directory"
# Two separate "x" insertions: one between "[ " and "]", another right after "]".
$start4 = Find-Start "[ ] Create Phase 2 working directory"
$posA = $start4 + 2   # between "[ " and "]"
$posB = $start4 + 4   # between "]" and " Create..." (after the first "x" shifts it by 1)
$d.Range($posA, $posA).InsertBefore("x")
$d.Range($posB, $posB).InsertBefore("x")
# Split right-to-left so each split's right-hand neighbour is already final.
Split-Char $posB
Split-Char $posA

# "[ ] Review and approve detailed plan"  ->  "[ x] Review and approve detailed plan"
Insert-CheckMark "[ ] Review and approve detailed plan" 2

# --- Week 1: YAML Enhancement --------------------------------------------
# "[ ] Add recognition_config section"  ->  "[x ] Add recognition_config section"
# (the "x" lands right after "[", before the existing space)
Insert-CheckMark "[ ] Add recognition_config section" 1

# "[ ] Add hospital_overrides section"  ->  "[ x] Add hospital_overrides section"
Insert-CheckMark "[ ] Add hospital_overrides section" 2

# "[ ] Template structure"  ->  "[ x] Template structure"
Insert-CheckMark "[ ] Template structure" 2

Write-Host "Done"
